# Remove the post entry that was row 417 ("「創造的であり続けるための40の方法」...").
# Excel's normal row-delete semantics: deletes the entire row and shifts
# every row below it up by one, which also updates the used range
# (dimension) from A1:C619 to A1:C618 — matching the target diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(417).Delete()
